# Issue #9: read in the health cost params and calculate total cost per infection
$wb = $excel.ActiveWorkbook

# Rename the "health-costs" sheet to "healthcosts"
$ws = $wb.Worksheets.Item("health-costs")
$ws.Name = "healthcosts"

# Add a "severity" header label in column A (new shared string)
$ws.Range("A1").Value = "severity"

# Turn the chronic DALY_weight value into an explicit formula (same literal value)
$ws.Range("B5").Formula = "=0.219"

# Widen the columns that will hold the new cost-per-infection calculations
$ws.Columns.Item(9).ColumnWidth = 29.5
$ws.Columns.Item(10).ColumnWidth = 29.333333333333332
$ws.Columns.Item(11).ColumnWidth = 29.333333333333332

# Update the active selection to reflect where work continued
$ws.Range("H27").Select() | Out-Null
